$wb = $excel.ActiveWorkbook

# --- Work on the "Тест" (Test) sheet: add 3 new test-unit rows ---
$wsTest = $wb.Worksheets.Item("Тест")
$wsTest.Activate()

# Unit names first (column A for the two "manekin" rows)
$wsTest.Range("A3").Value = "Яростный манекен"
$wsTest.Range("A4").Value = "Стойкий яростный манекен"

# Shared ability text for rows 3 & 4
$wsTest.Range("O3").Value = "Живое существо, Гнев манекена"
$wsTest.Range("O4").Value = "Живое существо, Гнев манекена"

# Ability text for row 5, then its unit name
$wsTest.Range("O5").Value = "Ослабляющий удар, Нежить"
$wsTest.Range("A5").Value = "Манекен зомби"

# Remaining numeric stats, row by row
$wsTest.Range("B3").Value = 0
$wsTest.Range("C3").Value = 0
$wsTest.Range("D3").Value = 0
$wsTest.Range("E3").Value = 0
$wsTest.Range("F3").Value = 100
$wsTest.Range("G3").Value = 1
$wsTest.Range("H3").Value = 0
$wsTest.Range("I3").Value = 0
$wsTest.Range("J3").Value = 0
$wsTest.Range("K3").Value = 1
$wsTest.Range("L3").Value = 1
$wsTest.Range("M3").Value = 1
$wsTest.Range("N3").Value = 0
$wsTest.Range("P3").Value = 1
$wsTest.Range("Q3").Value = 1
$wsTest.Range("R3").Value = 0

$wsTest.Range("B4").Value = 0
$wsTest.Range("C4").Value = 100
$wsTest.Range("D4").Value = 0
$wsTest.Range("E4").Value = 0
$wsTest.Range("F4").Value = 100
$wsTest.Range("G4").Value = 1
$wsTest.Range("H4").Value = 0
$wsTest.Range("I4").Value = 0
$wsTest.Range("J4").Value = 0
$wsTest.Range("K4").Value = 1
$wsTest.Range("L4").Value = 1
$wsTest.Range("M4").Value = 1
$wsTest.Range("N4").Value = 0
$wsTest.Range("P4").Value = 1
$wsTest.Range("Q4").Value = 1
$wsTest.Range("R4").Value = 0

$wsTest.Range("B5").Value = 0
$wsTest.Range("C5").Value = 100
$wsTest.Range("D5").Value = 0
$wsTest.Range("E5").Value = 0
$wsTest.Range("F5").Value = 100
$wsTest.Range("G5").Value = 1
$wsTest.Range("H5").Value = 0
$wsTest.Range("I5").Value = 0
$wsTest.Range("J5").Value = 0
$wsTest.Range("K5").Value = 1
$wsTest.Range("L5").Value = 1
$wsTest.Range("M5").Value = 1
$wsTest.Range("N5").Value = 0
$wsTest.Range("P5").Value = 1
$wsTest.Range("Q5").Value = 1
$wsTest.Range("R5").Value = 0

# Widen column A on the "Тест" sheet to fit the longer unit names
$wsTest.Columns.Item(1).ColumnWidth = 28.6

# Bugfix: move the selection/active cell on "Тест" to O4 and make it the tab shown when the workbook opens
$null = $wsTest.Range("O4").Select()

# --- The previously active sheet "Орден порядка" loses tab focus, selection moves to O16 ---
$wsOrder = $wb.Worksheets.Item("Орден порядка")
$wsOrder.Activate()
$null = $wsOrder.Range("O16").Select()

# Re-activate "Тест" so it ends up as the active/selected tab on save
$wsTest.Activate()
